$d = $word.ActiveDocument

function Add-PeriodAfter([string]$afterText) {
    # Locate the target text and remember the position right after it.
    $r = $d.Content
    $found = $r.Find.Execute($afterText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $afterText"
    }
    $pos = $r.End

    # Borrow formatting (Times New Roman, 12pt) from the pre-existing lone
    # "." run elsewhere in the document so the new run's rPr matches exactly.
    $refR = $d.Content
    $foundRef = $refR.Find.Execute("BillignLineItems.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $foundRef) {
        throw "Could not find reference run for formatting"
    }
    $dotRange = $d.Range($refR.End - 1, $refR.End)
    $ft = $dotRange.FormattedText

    # Insert a placeholder character right after the target text, then
    # stamp it with the borrowed formatted text ("." + matching rPr).
    $target = $d.Range($pos, $pos)
    $target.InsertAfter("X")
    $target2 = $d.Range($pos, $pos + 1)
    $target2.FormattedText = $ft

    # Return the range that now spans the newly-inserted "."
    return $d.Range($pos, $pos + 1)
}

# 1/2. Add the missing trailing "." for the Nurse and Doctor captions.
Add-PeriodAfter "ws how the Nurse can interact with PIMS" | Out-Null
$doctorDot = Add-PeriodAfter "ws how the Doctor can interact with PIMS"

# 3. Move the "_GoBack" bookmark from right after "Volunteer State Diagram"
#    to the very end of the document (after the Doctor caption's new ".").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Placing a bookmark exactly at the very end of the document's content can
# silently collapse to position 0, so temporarily extend the document with
# placeholder text, add the bookmark, then remove the placeholder again.
# The bookmark stays correctly anchored at the original end position.
$endPos = $d.Content.End
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("ZZPLACEHOLDERZZ")

$bookmarkRange = $d.Range($doctorDot.End, $doctorDot.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$phRange = $d.Content
$foundPh = $phRange.Find.Execute("ZZPLACEHOLDERZZ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPh) {
    $phRange.Text = ""
}
